# Scheduled-runner style update of recalculated market/profit figures
# (columns H:N) across several per-job sheets. Values below are the
# refreshed currentAveragePrice*/LevePrice*/LeveProfit* numbers.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(38, 8).Value = 123.14286
$ws.Cells.Item(38, 9).Value = 123.14286
$ws.Cells.Item(38, 11).Value = 369.42858
$ws.Cells.Item(38, 13).Value = 2.571419999999989

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 3895.8333
$ws.Cells.Item(62, 9).Value = 3888.8462
$ws.Cells.Item(62, 11).Value = 3888.8462
$ws.Cells.Item(62, 13).Value = -3264.8462

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(65, 8).Value = 3895.8333
$ws.Cells.Item(65, 9).Value = 3888.8462
$ws.Cells.Item(65, 11).Value = 19444.231
$ws.Cells.Item(65, 13).Value = -16324.231

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(86, 8).Value = 5342.143
$ws.Cells.Item(86, 9).Value = 8666.666999999999
$ws.Cells.Item(86, 10).Value = 2848.75
$ws.Cells.Item(86, 11).Value = 8666.666999999999
$ws.Cells.Item(86, 12).Value = 2848.75
$ws.Cells.Item(86, 13).Value = -7543.666999999999
$ws.Cells.Item(86, 14).Value = -5094.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(89, 8).Value = 5342.143
$ws.Cells.Item(89, 9).Value = 8666.666999999999
$ws.Cells.Item(89, 10).Value = 2848.75
$ws.Cells.Item(89, 11).Value = 43333.335
$ws.Cells.Item(89, 12).Value = 14243.75
$ws.Cells.Item(89, 13).Value = -37717.335
$ws.Cells.Item(89, 14).Value = -25475.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 3407.1428
$ws.Cells.Item(137, 9).Value = 3210.75
$ws.Cells.Item(137, 11).Value = 9632.25
$ws.Cells.Item(137, 13).Value = -7082.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2898.125
$ws.Cells.Item(2, 9).Value = 2897.5
$ws.Cells.Item(2, 11).Value = 2897.5
$ws.Cells.Item(2, 13).Value = -2784.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(116, 8).Value = 2898.125
$ws.Cells.Item(116, 9).Value = 2897.5
$ws.Cells.Item(116, 11).Value = 2897.5
$ws.Cells.Item(116, 13).Value = -603.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 3429.3333
$ws.Cells.Item(132, 9).Value = 3429.3333
$ws.Cells.Item(132, 11).Value = 10287.9999
$ws.Cells.Item(132, 13).Value = -7757.999899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2898.125
$ws.Cells.Item(3, 9).Value = 2897.5
$ws.Cells.Item(3, 11).Value = 2897.5
$ws.Cells.Item(3, 13).Value = -2783.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 1661.3158
$ws.Cells.Item(94, 9).Value = 1659.375
$ws.Cells.Item(94, 10).Value = 1671.6666
$ws.Cells.Item(94, 11).Value = 1659.375
$ws.Cells.Item(94, 12).Value = 1671.6666
$ws.Cells.Item(94, 13).Value = -1208.375
$ws.Cells.Item(94, 14).Value = -2573.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 793.06665
$ws.Cells.Item(107, 9).Value = 706.8570999999999
$ws.Cells.Item(107, 10).Value = 2000
$ws.Cells.Item(107, 11).Value = 706.8570999999999
$ws.Cells.Item(107, 12).Value = 2000
$ws.Cells.Item(107, 13).Value = 1213.1429
$ws.Cells.Item(107, 14).Value = -5840

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(137, 8).Value = 41666

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2649.3572
$ws.Cells.Item(31, 9).Value = 1924.1818
$ws.Cells.Item(31, 10).Value = 5308.3335
$ws.Cells.Item(31, 11).Value = 1924.1818
$ws.Cells.Item(31, 12).Value = 5308.3335
$ws.Cells.Item(31, 13).Value = -1629.1818
$ws.Cells.Item(31, 14).Value = -5898.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 2649.3572
$ws.Cells.Item(34, 9).Value = 1924.1818
$ws.Cells.Item(34, 10).Value = 5308.3335
$ws.Cells.Item(34, 11).Value = 1924.1818
$ws.Cells.Item(34, 12).Value = 5308.3335
$ws.Cells.Item(34, 13).Value = -1722.1818
$ws.Cells.Item(34, 14).Value = -5712.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 3465.3809
$ws.Cells.Item(58, 9).Value = 2617.4375
$ws.Cells.Item(58, 10).Value = 6178.8
$ws.Cells.Item(58, 11).Value = 2617.4375
$ws.Cells.Item(58, 12).Value = 6178.8
$ws.Cells.Item(58, 13).Value = -2414.4375
$ws.Cells.Item(58, 14).Value = -6584.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(94, 8).Value = 1817.8334
$ws.Cells.Item(94, 9).Value = 1841.6
$ws.Cells.Item(94, 11).Value = 1841.6
$ws.Cells.Item(94, 13).Value = -1390.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 1582.4166
$ws.Cells.Item(132, 9).Value = 1561.25
$ws.Cells.Item(132, 10).Value = 1624.75
$ws.Cells.Item(132, 11).Value = 4683.75
$ws.Cells.Item(132, 12).Value = 4874.25
$ws.Cells.Item(132, 13).Value = -2153.75
$ws.Cells.Item(132, 14).Value = -9934.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 3465.3809
$ws.Cells.Item(136, 9).Value = 2617.4375
$ws.Cells.Item(136, 10).Value = 6178.8
$ws.Cells.Item(136, 11).Value = 7852.3125
$ws.Cells.Item(136, 12).Value = 18536.4
$ws.Cells.Item(136, 13).Value = -5302.3125
$ws.Cells.Item(136, 14).Value = -23636.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 378.8
$ws.Cells.Item(2, 9).Value = 376.14285
$ws.Cells.Item(2, 11).Value = 2256.8571
$ws.Cells.Item(2, 13).Value = -2143.8571

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(14, 8).Value = 50099.5
$ws.Cells.Item(14, 9).Value = 50099.5
$ws.Cells.Item(14, 11).Value = 150298.5
$ws.Cells.Item(14, 13).Value = -150125.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(86, 8).Value = 1162
$ws.Cells.Item(86, 9).Value = 987
$ws.Cells.Item(86, 10).Value = 1249.5
$ws.Cells.Item(86, 11).Value = 2961
$ws.Cells.Item(86, 12).Value = 3748.5
$ws.Cells.Item(86, 13).Value = -1775
$ws.Cells.Item(86, 14).Value = -6120.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(89, 8).Value = 1162
$ws.Cells.Item(89, 9).Value = 987
$ws.Cells.Item(89, 10).Value = 1249.5
$ws.Cells.Item(89, 11).Value = 8883
$ws.Cells.Item(89, 12).Value = 11245.5
$ws.Cells.Item(89, 13).Value = -2955
$ws.Cells.Item(89, 14).Value = -23101.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(92, 8).Value = 425
$ws.Cells.Item(92, 9).Value = 425
$ws.Cells.Item(92, 10).Value = 425
$ws.Cells.Item(92, 11).Value = 1275
$ws.Cells.Item(92, 12).Value = 1275
$ws.Cells.Item(92, 13).Value = -27
$ws.Cells.Item(92, 14).Value = -3771

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 594.75
$ws.Cells.Item(131, 9).Value = 668.1667
$ws.Cells.Item(131, 10).Value = 374.5
$ws.Cells.Item(131, 11).Value = 2004.5001
$ws.Cells.Item(131, 12).Value = 1123.5
$ws.Cells.Item(131, 13).Value = 3035.4999
$ws.Cells.Item(131, 14).Value = -11203.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(100, 8).Value = 2381.125
$ws.Cells.Item(100, 9).Value = 809.8
$ws.Cells.Item(100, 11).Value = 809.8
$ws.Cells.Item(100, 13).Value = -268.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 5168
$ws.Cells.Item(132, 9).Value = 336
$ws.Cells.Item(132, 10).Value = 10000
$ws.Cells.Item(132, 11).Value = 1008
$ws.Cells.Item(132, 12).Value = 30000
$ws.Cells.Item(132, 13).Value = 1522
$ws.Cells.Item(132, 14).Value = -35060

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 910.5833
$ws.Cells.Item(100, 9).Value = 994.44446
$ws.Cells.Item(100, 11).Value = 1988.88892
$ws.Cells.Item(100, 13).Value = -1447.88892

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 3148.3333
$ws.Cells.Item(122, 9).Value = 3120.7856
$ws.Cells.Item(122, 11).Value = 9362.356800000001
$ws.Cells.Item(122, 13).Value = -6912.356800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 3779.919
$ws.Cells.Item(136, 9).Value = 3626.2666
$ws.Cells.Item(136, 10).Value = 4438.4287
$ws.Cells.Item(136, 11).Value = 10878.7998
$ws.Cells.Item(136, 12).Value = 13315.2861
$ws.Cells.Item(136, 13).Value = -8328.799800000001
$ws.Cells.Item(136, 14).Value = -18415.2861
